$wb = $excel.ActiveWorkbook

# "Erläuterungen" sheet holds the term/definition table
$ws = $wb.Worksheets.Item("Erläuterungen")

# Update the definition for "Anzahl Wohnungen in Sample 1 bzw. Sample 2" (row 9, col B):
# remove the trailing "Methodendokumentation unter: ..." sentence
$ws.Range("B9").Value = "Samplegrösse für die betreffende Zelle pro Schicht 1 resp. Schicht 2: Anzahl Mietpreisinformationen, die vorliegen."

# Update the definition for "Abfragetool MPE" (row 10, col B): new URL
$ws.Range("B10").Value = "Die Detaildaten der Mietpreiserhebung 2022 sind auf einem Abruftool verfügbar, das erreichbar ist unter: https://www.stadt-zuerich.ch/prd/de/index/statistik/publikationen-angebote/datenbanken-anwendungen/mietpreiserhebung.html"

# Move selection to B10 on the active sheet (matches diff's sheetView selection change)
$ws.Range("B10").Select()
